$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K data for year 2022
$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 3.9462868231169921
$ws.Range("K6").Value = 3.8007658934388928

# K4: like J4's style but without the fill flag
$k4 = $ws.Range("K4")
$k4.Font.Name = "Times New Roman"
$k4.Font.Bold = $true
$k4.Font.Size = 9
$k4.Borders.Item(3).LineStyle = 1
$k4.Borders.Item(3).Weight = -4138
$k4.Borders.Item(4).LineStyle = 1
$k4.Borders.Item(4).Weight = -4138
$k4.HorizontalAlignment = -4152
$k4.VerticalAlignment = -4160
$k4.WrapText = $true

# K5: numFmt 164, font "Times New Roman Cyr" sz 9, right aligned, wrap text
$k5 = $ws.Range("K5")
$k5.NumberFormat = "0.0"
$k5.Font.Name = "Times New Roman Cyr"
$k5.Font.Size = 9
$k5.HorizontalAlignment = -4152
$k5.WrapText = $true

# K6: same as K5 but with a bottom medium border
$k6 = $ws.Range("K6")
$k6.NumberFormat = "0.0"
$k6.Font.Name = "Times New Roman Cyr"
$k6.Font.Size = 9
$k6.HorizontalAlignment = -4152
$k6.WrapText = $true
$k6.Borders.Item(4).LineStyle = 1
$k6.Borders.Item(4).Weight = -4138

# Update selection to reflect new active cell
$ws.Range("L5").Select() | Out-Null
